$d = $word.ActiveDocument

# Locate the paragraph that ends the lead-in sentence right before the
# forecasts table ("... summarized in the table below:").
$targetIndex = -1
$i = 0
foreach ($para in $d.Paragraphs) {
    $i = $i + 1
    if ($para.Range.Text -like "*summarized in the table below:*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq -1) {
    throw "Could not find the target paragraph."
}

$targetPara = $d.Paragraphs.Item($targetIndex)

# Insert a new, empty paragraph right after it.
$targetPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Style = "Source Code"
$newRange = $newPara.Range
$insertStart = $newRange.Start

# Use InsertXML (flat-OPC package fragment) so the two text runs come in
# with xml:space="preserve" on <w:t>, matching the rest of the document.
$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
  '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
    '<pkg:xmlData>' + `
      '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + `
          '<w:p>' + `
            '<w:pPr><w:pStyle w:val="SourceCode"/></w:pPr>' + `
            '<w:r><w:t xml:space="preserve">pander</w:t></w:r>' + `
            '<w:r><w:t xml:space="preserve">(table_forecasts)</w:t></w:r>' + `
          '</w:p>' + `
        '</w:body>' + `
      '</w:document>' + `
    '</pkg:xmlData>' + `
  '</pkg:part>' + `
'</pkg:package>'
$newRange.InsertXML($xml)

# Reapply the paragraph style (InsertXML replaces the whole paragraph) and
# restore the character styles on each run — rStyle doesn't survive
# InsertXML, but Range.Style assignment applies it correctly.
$newPara2 = $d.Paragraphs.Item($targetIndex + 1)
$newPara2.Style = "Source Code"

$functionRange = $d.Range($insertStart, $insertStart + 6)
$functionRange.Style = "FunctionTok"

$normalRange = $d.Range($insertStart + 6, $insertStart + 24)
$normalRange.Style = "NormalTok"
